$d = $word.ActiveDocument

# Row 1, Col 1: "Rafid" -> cleared
$d.Content.Find.Execute("Rafid", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# Row 1, Col 2: "1722006" -> "I playing the game"
$d.Content.Find.Execute("1722006", $true, $false, $false, $false, $false, $true, 1, $false, "I playing the game", 2)

# Row 2, Col 1: "SK Zaman" -> cleared
$d.Content.Find.Execute("SK Zaman", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# Row 2, Col 2: "1724568" -> cleared
$d.Content.Find.Execute("1724568", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# Row 3, Col 1: "Asif" -> cleared
$d.Content.Find.Execute("Asif", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# Row 3, Col 2: "1345678" -> cleared
$d.Content.Find.Execute("1345678", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
